$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark these rows' status column (D) as "As-is"
$ws.Range("D3").Value = "As-is"
$ws.Range("D6").Value = "As-is"
$ws.Range("D7").Value = "As-is"
$ws.Range("D11").Value = "As-is"
$ws.Range("D15").Value = "As-is"

# Update the visible selection/view to D16
$ws.Range("D16").Select()
